# Repull data, push all data, mean calculation
# Update the dSF (column F) values for several rows to reflect the repulled data.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("F4").Value = -2
$ws.Range("F5").Value = 3
$ws.Range("F9").Value = -1
$ws.Range("F12").Value = 3
$ws.Range("F13").Value = 7
$ws.Range("F18").Value = -4
$ws.Range("F21").Value = -1
$ws.Range("F24").Value = 3
